# The workbook tracks weekly wholesale "Arveja Verde" (green pea) prices.
# This commit adds one new weekly observation. It is inserted as a new
# row 118, which pushes all the existing rows 118-130 down by one
# (becoming rows 119-131), growing the sheet from A1:R130 to A1:R131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 118, shifting rows 118:130
# down to 119:131 (and carrying their formatting/styles along).
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new observation's data.
$ws.Cells.Item(118, 1).Value  = 6
$ws.Cells.Item(118, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(118, 3).Value  = "Metropolitana"
$ws.Cells.Item(118, 4).Value  = 44461
$ws.Cells.Item(118, 5).Value  = 13
$ws.Cells.Item(118, 6).Value  = 100112022
$ws.Cells.Item(118, 7).Value  = "Arveja Verde"
$ws.Cells.Item(118, 8).Value  = "Sin especificar"
$ws.Cells.Item(118, 9).Value  = "Primera"
$ws.Cells.Item(118, 10).Value = 250
$ws.Cells.Item(118, 11).Value = 25000
$ws.Cells.Item(118, 12).Value = 27000
$ws.Cells.Item(118, 13).Value = 26200
$ws.Cells.Item(118, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(118, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(118, 16).Value = 1048
$ws.Cells.Item(118, 17).Value = 25
$ws.Cells.Item(118, 18).Value = "Hortaliza"
